$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "-"

# Row 3
$ws.Range("C3").Value = "[-, -, 'MEC-3B-Metrologia 2', -]"
$ws.Range("D3").Value = "-"

# Row 4
$ws.Range("C4").Value = "[-, -, 'MEC-3B-Metrologia 2', -]"
$ws.Range("D4").Value = "-"

# Row 6
$ws.Range("C6").Value = "[-, -, 'MEC-3B-Metrologia 2', -]"
$ws.Range("D6").Value = "-"

# Row 7
$ws.Range("C7").Value = "[-, -, 'MEC-3B-Metrologia 2', -]"

# Row 11
$ws.Range("C11").Value = "[-, 'MEC-1A-Metrologia 1', -, -]"

# Row 12
$ws.Range("C12").Value = "[-, 'MEC-1A-Metrologia 1', -, -]"

# Row 14
$ws.Range("C14").Value = "[-, 'MEC-1A-Metrologia 1', -, -]"

# Row 15
$ws.Range("C15").Value = "[-, 'MEC-1A-Metrologia 1', -, -]"

# Row 18
$ws.Range("C18").Value = "[-, 'MEC-2NB-Metrologia 2', -, -]"
$ws.Range("D18").Value = "[-, 'MEC-2NA-Metrologia 2', -, -]"

# Row 19
$ws.Range("C19").Value = "[-, 'MEC-2NB-Metrologia 2', -, -]"
$ws.Range("D19").Value = "[-, 'MEC-2NA-Metrologia 2', -, -]"

# Row 20
$ws.Range("C20").Value = "[-, 'MEC-2NB-Metrologia 2', -, -]"
$ws.Range("D20").Value = "[-, 'MEC-2NA-Metrologia 2', -, -]"
$ws.Range("E20").Value = "[-, 'MEC-2NB-Metrologia 2', -, -]"

# Row 21
$ws.Range("C21").Value = "-"
$ws.Range("D21").Value = "[-, 'MEC-2NA-Metrologia 2', -, -]"
